# Applies the edit described by the diff:
#  - Remove the stray row (old row 13) whose only content was the teacher's
#    name sitting in columns B/C with no label in column A. Deleting this
#    row shifts everything below it up by one.
#  - The text content also gets rearranged: several long paragraphs are
#    replaced by shorter text (and one formerly unlabeled value migrates to
#    become the visible value for "Objetivos:"/"Método:"). Set each affected
#    cell to its corrected text explicitly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the stray row (former row 13: only B13/C13 populated with the
# teacher name, no label in column A). This shifts rows 14-25 up to 13-24.
$ws.Rows(13).Delete()

# Fix up the values that are now misaligned after the shift / rewrite.
$ws.Range("B10").Value = "8767640 - Eduardo Ferro dos Santos"
$ws.Range("C10").Value = "8767640 - Eduardo Ferro dos Santos"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("B15").Value = "01/01/2021"
$ws.Range("C15").Value = "01/01/2021"

$ws.Range("B18").Value = "8767640 - Eduardo Ferro dos Santos"
$ws.Range("C18").Value = "8767640 - Eduardo Ferro dos Santos"

$ws.Range("B19").Value = "Aulas expositivas e práticas."
$ws.Range("C19").Value = "Aulas expositivas e práticas."

$ws.Range("B20").Value = "Exercícios de aprendizado e exercícios de avaliação farão parte da composição de notas individuais (NI), com aplicação de trabalhos práticos em grupo (NG). Sendo: Nota Final = (NI+NG)/2"
$ws.Range("C20").Value = "Exercícios de aprendizado e exercícios de avaliação farão parte da composição de notas individuais (NI), com aplicação de trabalhos práticos em grupo (NG). Sendo: Nota Final = (NI+NG)/2"

$ws.Range("B21").Value = "A recuperação deverá consistir de uma prova englobando a matéria toda do semestre. - A média final (pós-recuperação) deverá ser composta por uma média simples entre a nota do semestre (nota final) e a da prova de recuperação."
$ws.Range("C21").Value = "A recuperação deverá consistir de uma prova englobando a matéria toda do semestre. - A média final (pós-recuperação) deverá ser composta por uma média simples entre a nota do semestre (nota final) e a da prova de recuperação."
